# "Reports - Final version"
# Restructure the workbook: insert a new "Template_M" worksheet (a condensed
# report-template picker) ahead of the existing "Template" sheet, and extend
# "Template" itself with new Schedule/Frequency columns while trimming it
# down to a single data row. "Users" is left untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create "Template_M" by copying "Template" (placed immediately before
#    it), then rename the copy. NOTE: after Copy(), any previously-held
#    worksheet object reference tracks the *new* copy (not the original
#    sheet it was copied from) -- so sheets must be re-fetched by name
#    after this call.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Template").Copy($wb.Worksheets.Item("Template"))
$wb.Worksheets.Item("Template (2)").Name = "Template_M"

$templateM = $wb.Worksheets.Item("Template_M")
$template = $wb.Worksheets.Item("Template")

# ---------------------------------------------------------------------
# 2. Reshape "Template_M": drop the old "Account Owner Name" / "Shadow
#    account name" columns (L:M) and insert one new column after "Access"
#    (column C) for the View/Public/Private flag, then rewrite all values.
# ---------------------------------------------------------------------
$templateM.Columns("L:M").Delete()
$templateM.Columns.Item(3).Insert()

$templateM.Cells.ClearContents()

# Row 1 - headers
$templateM.Range("A1").Value = "Report Type"
$templateM.Range("B1").Value = "Template Name"
$templateM.Range("C1").Value = "Access"
$templateM.Range("D1").Value = "Currency"
$templateM.Range("E1").Value = "Amount Range"
$templateM.Range("F1").Value = "View"
$templateM.Range("G1").Value = "Account type"
$templateM.Range("H1").Value = "Transaction Type"
$templateM.Range("I1").Value = "Period"
$templateM.Range("J1").Value = "Optional fields"
$templateM.Range("K1").Value = "Include currency accounts"

# Row 2
$templateM.Range("A2").Value = "Transaction"
$templateM.Range("B2").Value = "Trans_"
$templateM.Range("C2").Value = "Public"
$templateM.Range("D2").Value = "Select All"
$templateM.Range("E2").Value = "-99999 - 99999999"
$templateM.Range("F2").Value = "Reporting Level"
$templateM.Range("G2").Value = "Virtual Transaction Account"
$templateM.Range("H2").Value = "All"
$templateM.Range("I2").Value = "Current month"
$templateM.Range("J2").Value = "Account name|Account Owner customer ID|Account Owner Name|Shadow account number|Shadow account name"
$templateM.Range("K2").Value = "Yes"

# Row 3
$templateM.Range("A3").Value = "Transaction"
$templateM.Range("B3").Value = "Trans_"
$templateM.Range("C3").Value = "Private"
$templateM.Range("D3").Value = "NOK"
$templateM.Range("E3").Value = "-99999 - 99999999"
$templateM.Range("G3").Value = "Virtual Transaction Account"
$templateM.Range("H3").Value = "All"
$templateM.Range("I3").Value = "Previous month"
$templateM.Range("J3").Value = "Account name|Account Owner customer ID|Account Owner Name|Shadow account number|Shadow account name"

# ---------------------------------------------------------------------
# 3. Reshape "Template": keep only the first data row, and add the new
#    "Select Schedule" / "Frequency" columns (L:M).
# ---------------------------------------------------------------------
$template.Rows.Item(3).Delete()

$template.Range("L1").Value = "Select Schedule"
$template.Range("M1").Value = "Frequency"
$template.Range("L2").Value = "Create Schedule to Run Later"
$template.Range("M2").Value = "Every Week"

$template.Range("F1").Value = "View"
$template.Range("J1").Value = "Optional fields"
$template.Range("K1").Value = "Include currency accounts"
$template.Range("F2").Value = "Reporting Level"
$template.Range("H2").Value = "All"
$template.Range("J2").Value = "Account name|Account Owner customer ID|Account Owner Name|Shadow account number|Shadow account name"

# ---------------------------------------------------------------------
# 4. View state: "Template" is the active/selected tab, scrolled toward
#    the new columns; "Template_M" keeps a plain view with row 3 selected.
# ---------------------------------------------------------------------
$templateM = $wb.Worksheets.Item("Template_M")
$templateM.Range("A3:XFD3").Select()

$template = $wb.Worksheets.Item("Template")
$template.Activate()
$template.Range("L1").Select()
